$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the neighboring "sum" header cell (G1) so the
# new "Save" header cell (H1) gets the same bold/centered/bordered format.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the new "Save" column: header in row 1, value in row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
